$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- Row 95: headers ---
$ws.Range("S95").Value = "outside-to-outside outer hole to diag outer hole"
$ws.Range("W95").Value = "outside-to-outside outer hole to adjacent outer hole"
$ws.Range("AC95").Value = "hole diameter"

# --- Rows 96-99: raw caliper measurements (numFmt 0.00) ---
$ws.Range("S96").Value = 110.66
$ws.Range("S96").NumberFormat = "0.00"
$ws.Range("W96").Value = 85.62
$ws.Range("W96").NumberFormat = "0.00"
$ws.Range("Y96").Value = 70.569999999999993
$ws.Range("Y96").NumberFormat = "0.00"
$ws.Range("AC96").Value = 1.45

$ws.Range("S97").Value = 110.45
$ws.Range("S97").NumberFormat = "0.00"
$ws.Range("W97").Value = 86.01
$ws.Range("W97").NumberFormat = "0.00"
$ws.Range("Y97").Value = 70.83
$ws.Range("Y97").NumberFormat = "0.00"
$ws.Range("AC97").Value = 1.54

$ws.Range("S98").Value = 110.61
$ws.Range("S98").NumberFormat = "0.00"
$ws.Range("W98").Value = 85.88
$ws.Range("W98").NumberFormat = "0.00"
$ws.Range("Y98").Value = 70.61
$ws.Range("Y98").NumberFormat = "0.00"
$ws.Range("AC98").Value = 1.67

$ws.Range("S99").Value = 110.7
$ws.Range("S99").NumberFormat = "0.00"
$ws.Range("W99").Value = 85.6
$ws.Range("W99").NumberFormat = "0.00"
$ws.Range("Y99").Value = 70.69
$ws.Range("Y99").NumberFormat = "0.00"
$ws.Range("AC99").Value = 1.43

# --- Row 100: averages (bold, 0.00) + distance (8 decimals) ---
$ws.Range("S100").Formula = "=AVERAGE(S96:S99)"
$ws.Range("S100").NumberFormat = "0.00"
$ws.Range("S100").Font.Bold = $true

$ws.Range("W100").Formula = "=AVERAGE(W96:W99)"
$ws.Range("W100").NumberFormat = "0.00"
$ws.Range("W100").Font.Bold = $true

$ws.Range("Y100").Formula = "=AVERAGE(Y96:Y99)"
$ws.Range("Y100").NumberFormat = "0.00"
$ws.Range("Y100").Font.Bold = $true

$ws.Range("AA100").Formula = "=SQRT(POWER(W100,2)+POWER(Y100,2))"
$ws.Range("AA100").NumberFormat = "0.00000000"

$ws.Range("AC100").Formula = "=AVERAGE(AC96:AC99)"
$ws.Range("AC100").NumberFormat = "0.00"
$ws.Range("AC100").Font.Bold = $true

# --- Row 101: carries styles only, no values ---
$ws.Range("AA101").NumberFormat = "0.00000000"
$ws.Range("AC101").NumberFormat = "0.00"
$ws.Range("AC101").Font.Bold = $true

# --- Row 102 ---
$ws.Range("S102").Formula = '=S100-0.5*$AC$100'
$ws.Range("W102").Formula = '=W100-0.5*$AC$100'
$ws.Range("Y102").Formula = '=Y100-0.5*$AC$100'
$ws.Range("AA102").Formula = "=SQRT(POWER(W102,2)+POWER(Y102,2))"
$ws.Range("AA102").NumberFormat = "0.00000000"

Write-Output "done-part1"
